$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells
$ws.Range("B2").Value = 194
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 104

# Remove rows 4 and 5 entirely (last cases dropped)
$ws.Range("A4:B5").ClearContents()
$ws.Rows("4:5").Delete()
